# "Generate Report for Handoff"
#
# A fresh handoff/XLIFF-generation run happened: the source file got a new
# GUID-based name, fresh handoff files + timestamps were produced, and
# (since nothing has been handed *back* yet for this new handoff) the
# "Latest Target/Handback" columns on the locale sheets are reset to empty /
# the zero-date sentinel.

$wb = $excel.ActiveWorkbook

$newGuidFile      = "08f44789-0b50-403a-9333-29418f4f5cfa.md"
$newGuidPath      = "e2e\08f44789-0b50-403a-9333-29418f4f5cfa.md"

$newHoDate        = "2016-08-28 17:00:24"   # Overview!G2 and de-de!H2 (shared text)
$newZhHoDate      = "2016-08-28 17:00:19"   # zh-cn!H2
$resetHandback    = "0001-01-01 00:00:00"   # zh-cn!K2 and de-de!K2

$newZhXlf = "08f44789-0b50-403a-9333-29418f4f5cfa.313345006f1057710f4080f26eae94fc9a9c4345.zh-cn.xlf"
$newDeXlf = "08f44789-0b50-403a-9333-29418f4f5cfa.313345006f1057710f4080f26eae94fc9a9c4345.de-de.xlf"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newGuidFile
$wsOverview.Range("B2").Value = $newGuidPath
$wsOverview.Range("G2").Value = $newHoDate

foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.Range.Address() -eq '$B$2') {
        $hl.TextToDisplay = $newGuidPath
    }
}

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newGuidFile
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = $newZhHoDate
$wsZh.Range("K2").Value = $resetHandback

foreach ($hl in $wsZh.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $newGuidFile
    } elseif ($addr -eq '$I$2') {
        $hl.Delete()
    }
}

# Latest Target File / Latest Handback File: no handback yet against the
# freshly generated handoff, so these clear out.
$wsZh.Range("I2").Value = ""
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("J2").Value = ""

$wsZh.Columns.Item(9).ColumnWidth = 17.85
$wsZh.Columns.Item(10).ColumnWidth = 20.85

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newGuidFile
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = $newHoDate
$wsDe.Range("K2").Value = $resetHandback

foreach ($hl in $wsDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $newGuidFile
    } elseif ($addr -eq '$I$2') {
        $hl.Delete()
    }
}

$wsDe.Range("I2").Value = ""
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("J2").Value = ""

$wsDe.Columns.Item(9).ColumnWidth = 17.85
$wsDe.Columns.Item(10).ColumnWidth = 20.85
